$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$ttle  = $wb.Worksheets.Item("TTLE")

# --- Content edits on the "About" sheet --------------------------------
# Replace the old one-line note with the new wording.
$about.Range("B3").Value = "None needed.  Handled through calibration."

# Add a "Notes" explanation block under the existing note.
$about.Range("A5").Value = "Notes"
$about.Range("A5").Font.Bold = $true

$about.Range("A6").Value  = "The logit exponents express how large of a cost difference between technology options"
$about.Range("A7").Value  = "is required to produce a change in technology selection.  This parameter needs to be"
$about.Range("A8").Value  = "obtained via model calibration - e.g. testing a given price intervention with different"
$about.Range("A9").Value  = "logit exponent values until it produces a technology choice shift that matches real-world"
$about.Range("A10").Value = "data on technology buyers' behavior."

$about.Range("A12").Value = 'For more on this, see the "Modified Logit" equation description at:'
$about.Range("A13").Value = "https://jgcri.github.io/gcam-doc/choice.html"

# --- Sheet activation / selection ---------------------------------------
# The saved workbook now opens on the "About" sheet (instead of "TTLE"),
# with both sheets' selections reset back to their top-left cell.
$ttle.Activate()
$ttle.Range("A1").Select()

$about.Activate()
$about.Range("A1").Select()

Write-Output "done"
